$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),`n                ('model', SVC(C=5, class_weight='balanced', random_state=42))])"
$ws.Range("B2").Value = 0.6857142857142857
$ws.Range("C2").Value = "{'scaler': None, 'model__kernel': 'rbf', 'model__class_weight': 'balanced', 'model__C': 5}"
$ws.Range("D2").Value = 0.4285714285714285
$ws.Range("E2").Value = "[1 0 0 1 0 0 1 1 0 1 0 0]"
$ws.Range("F2").Value = "[0 0 1 0 1 1 1 1 1 1 1 1]"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.6792559523809524
$ws.Range("I2").Value = 0.0354950218135263
$ws.Range("J2").Value = 0.5683333333333332
$ws.Range("K2").Value = 0.06457409615698313

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),`n                ('model', SVC(C=1, class_weight='balanced', random_state=42))])"
$ws.Range("B3").Value = 0.6571428571428571
$ws.Range("C3").Value = "{'scaler': None, 'model__kernel': 'rbf', 'model__class_weight': 'balanced', 'model__C': 1}"
$ws.Range("D3").Value = 0.3636363636363636
$ws.Range("E3").Value = "[1 0 1 0 0 0 0 1 1 0 1 1]"
$ws.Range("F3").Value = "[1 1 0 1 1 0 0 0 0 0 1 0]"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.6698511904761906
$ws.Range("I3").Value = 0.04027462134673227
$ws.Range("J3").Value = 0.5591071428571428
$ws.Range("K3").Value = 0.06360185880433948

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),`n                ('model', SVC(C=0.0001, kernel='linear', random_state=42))])"
$ws.Range("B4").Value = 0.6380952380952382
$ws.Range("C4").Value = "{'scaler': None, 'model__kernel': 'linear', 'model__class_weight': None, 'model__C': 0.0001}"
$ws.Range("D4").Value = 0.7058823529411765
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 0 1 1 1 1 1 0 1 0 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6598065476190476
$ws.Range("I4").Value = 0.02963008882387964
$ws.Range("J4").Value = 0.536845238095238
$ws.Range("K4").Value = 0.06408047362222596
